$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "alcohol data" sheet has an extra column (M) that duplicates/overlaps
# the next one; remove it so the following column (old N) shifts left to
# become the new M, matching the author's "remove column from alcohol
# data" change. Column 13 = M.
$ws.Columns.Item(13).Delete()

$ws.Activate()
$ws.Range("M1").Select()
